# Generate Report for handback
#
# The localization-status workbook tracks, per language, the handoff /
# handback lifecycle of each source file. This run marks the single
# tracked file (48a100d7-...md) as handed back (in sync with en-US) for
# both locales, records the target/handback file links, and stamps the
# handback datetime.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

$mdName  = "48a100d7-6b2b-40a1-8923-ea1b68988aeb.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/7f2e3581747e61c130943b9046f6f5b2d117e79a/e2e/48a100d7-6b2b-40a1-8923-ea1b68988aeb.md"

$zhXlfName = "48a100d7-6b2b-40a1-8923-ea1b68988aeb.072bcdca262c36e657437155c39331959b9ea964.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0a2c508075875462eb281961673e7f4ae88b693/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/48a100d7-6b2b-40a1-8923-ea1b68988aeb.072bcdca262c36e657437155c39331959b9ea964.zh-cn.xlf"

$deXlfName = "48a100d7-6b2b-40a1-8923-ea1b68988aeb.072bcdca262c36e657437155c39331959b9ea964.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70f0ebef98e6bb39c9826f77fe0ce10b4ed1a825/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/48a100d7-6b2b-40a1-8923-ea1b68988aeb.072bcdca262c36e657437155c39331959b9ea964.de-de.xlf"

# --- Overview sheet: roll the new status up for both locale columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $status
$overview.Range("C2").Value = $status

# --- zh-cn sheet: handback complete ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $status
$zh.Hyperlinks.Add($zh.Range("E2"), $mdUrl, $null, $null, $mdName)
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfUrl, $null, $null, $zhXlfName)
$zh.Range("G2").Value = "2016-01-17 14:40:49"

# --- de-de sheet: handback complete ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $status
$de.Hyperlinks.Add($de.Range("E2"), $mdUrl, $null, $null, $mdName)
$de.Hyperlinks.Add($de.Range("F2"), $deXlfUrl, $null, $null, $deXlfName)
$de.Range("G2").Value = "2016-01-17 14:41:08"
